# Applies the "Anzahl Kinder" -> "Anzahl Kinder unter 25" update plus the
# refreshed PV threshold figures and entry date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Anzahl Kinder unter 25"
$ws.Range("B4").Value = 62100
$ws.Range("B5").Value = 69300
$ws.Range("B6").Value = "01.01.2024"

$ws.Range("B9").Select()
